# Applies the "output generated at 456a3b4" refresh to 广州-漫展信息.xlsx
# Updates the "想去人数" (F) / "最低票价" (G) columns across the four sheets
# (展览, 演出, 本地生活, 全部类型) with newly scraped numbers, and marks two
# rows (无法订票) with G = "不可售" instead of a numeric minimum price.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G2").Value = "不可售"
$ws.Range("F3").Value = 555
$ws.Range("F4").Value = 239
$ws.Range("F5").Value = 19
$ws.Range("F6").Value = 717
$ws.Range("F7").Value = 333
$ws.Range("F9").Value = 125
$ws.Range("F10").Value = 235
$ws.Range("F11").Value = 194
$ws.Range("G11").Value = 55
$ws.Range("F12").Value = 5115
$ws.Range("F14").Value = 27
$ws.Range("F15").Value = 477
$ws.Range("F17").Value = 529
$ws.Range("F18").Value = 324
$ws.Range("F19").Value = 412
$ws.Range("F21").Value = 13
$ws.Range("F23").Value = 82
$ws.Range("F24").Value = 297
$ws.Range("F25").Value = 992
$ws.Range("F27").Value = 1701
$ws.Range("F28").Value = 422
$ws.Range("F29").Value = 30

# --- Sheet 2: 演出 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 42
$ws.Range("F8").Value = 291

# --- Sheet 3: 本地生活 ------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 167

# --- Sheet 4: 全部类型 (aggregated view of the above three sheets) --------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G3").Value = "不可售"
$ws.Range("F6").Value = 167
$ws.Range("F7").Value = 555
$ws.Range("F8").Value = 239
$ws.Range("F9").Value = 19
$ws.Range("F10").Value = 717
$ws.Range("F12").Value = 333
$ws.Range("F14").Value = 125
$ws.Range("F15").Value = 235
$ws.Range("F16").Value = 194
$ws.Range("G16").Value = 55
$ws.Range("F17").Value = 5115
$ws.Range("F19").Value = 27
$ws.Range("F21").Value = 477
$ws.Range("F23").Value = 529
$ws.Range("F24").Value = 324
$ws.Range("F25").Value = 412
$ws.Range("F26").Value = 42
$ws.Range("F28").Value = 13
$ws.Range("F30").Value = 291
$ws.Range("F36").Value = 82
$ws.Range("F37").Value = 297
$ws.Range("F38").Value = 992
$ws.Range("F40").Value = 1701
$ws.Range("F41").Value = 422
$ws.Range("F42").Value = 30
